# Break out stock.yaml "backup" column: add a new column R ("backup") with
# default 0 for all existing rows, fix a stray Q52 value, and append six new
# monthly rows (2024-07-01 .. 2024-12-01) pulled in by the updated fetch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell R1 = "backup", formatted like the other headers ---
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R1").Value = "backup"

# --- 2. New column R populated with 0 for every existing data row (2-265) ---
for ($r = 2; $r -le 265; $r++) {
    $ws.Cells.Item($r, 18).Value = 0
}

# --- 3. Data correction: Q52 was 2, now 0 ---
$ws.Range("Q52").Value = 0

# --- 4. Copy date formatting from A265 down to the new date cells A266:A271 ---
$ws.Range("A265").Copy()
$ws.Range("A266:A271").PasteSpecial(-4122)  # xlPasteFormats

# --- 5. Append the six new monthly rows (266-271) ---
$newRows = @(
    @{ Row=266; A=45474; B=2655.99350762891;   C=2866.73372461043;  D=2651.859491469526; E=2766.120361328125; G=20457458; H=2024; I=7;  J=1; K=0; L=0; M=0; N=27; O=1; P=0; Q=0 },
    @{ Row=267; A=45505; B=2771.300291347343;  C=2777.974594250171; D=2495.460397205453; E=2688.5185546875;   G=21059904; H=2024; I=8;  J=1; K=0; L=0; M=0; N=31; O=0; P=0; Q=0 },
    @{ Row=268; A=45536; B=2706.75;            C=2824;               D=2583.949951171875; E=2795.550048828125; G=16680276; H=2024; I=9;  J=1; K=0; L=0; M=0; N=35; O=0; P=0; Q=0 },
    @{ Row=269; A=45566; B=2794.050048828125;  C=2833;               D=2594.300048828125; E=2695.85009765625;  G=12231593; H=2024; I=10; J=1; K=0; L=0; M=0; N=40; O=0; P=0; Q=0 },
    @{ Row=270; A=45597; B=2719;               C=2719;               D=2472.050048828125; E=2606.25;           G=12495812; H=2024; I=11; J=1; K=0; L=0; M=0; N=44; O=0; P=0; Q=2 },
    @{ Row=271; A=45627; B=2605.85009765625;   C=2733.75;            D=2436;              E=2443.5;            G=10684669; H=2024; I=12; J=1; K=0; L=0; M=0; N=48; O=0; P=0; Q=0 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value  = $row.A   # A Datetime
    $ws.Cells.Item($r, 2).Value  = $row.B   # B Open
    $ws.Cells.Item($r, 3).Value  = $row.C   # C High
    $ws.Cells.Item($r, 4).Value  = $row.D   # D Low
    $ws.Cells.Item($r, 5).Value  = $row.E   # E Close
    # F (Adj Close) left blank - no data for these rows
    $ws.Cells.Item($r, 7).Value  = $row.G   # G Volume
    $ws.Cells.Item($r, 8).Value  = $row.H   # H Year
    $ws.Cells.Item($r, 9).Value  = $row.I   # I Month
    $ws.Cells.Item($r, 10).Value = $row.J   # J Day
    $ws.Cells.Item($r, 11).Value = $row.K   # K Hour
    $ws.Cells.Item($r, 12).Value = $row.L   # L Minute
    $ws.Cells.Item($r, 13).Value = $row.M   # M Second
    $ws.Cells.Item($r, 14).Value = $row.N   # N Week
    $ws.Cells.Item($r, 15).Value = $row.O   # O isPivot
    $ws.Cells.Item($r, 16).Value = $row.P   # P two_line_structure
    $ws.Cells.Item($r, 17).Value = $row.Q   # Q detect_structure
    # R (backup) left blank - no data for these rows
}
